$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map, derived from the authoritative diff.
# Every target cell in this sheet is stored as text (coinranking scraper
# writes everything as inline strings), including numeric-looking values
# like prices ("248.89") and hour-of-day ("12"), so each write below forces
# the cell to Text (NumberFormat "@") before assigning, then clears the
# number-format override again so no stray style is left behind on the cell
# (matches the unstyled/default-style cells in the source workbook).
$updates = @(
    @{ Cell = 'D2'; Value = '248.89' }
    @{ Cell = 'G2'; Value = '12' }
    @{ Cell = 'D3'; Value = '22.72' }
    @{ Cell = 'G3'; Value = '12' }
    @{ Cell = 'D4'; Value = '5.307' }
    @{ Cell = 'G4'; Value = '12' }
    @{ Cell = 'D5'; Value = '0.05695' }
    @{ Cell = 'G5'; Value = '12' }
    @{ Cell = 'D6'; Value = '3.407' }
    @{ Cell = 'G6'; Value = '12' }
    @{ Cell = 'D7'; Value = '6.347' }
    @{ Cell = 'G7'; Value = '12' }
    @{ Cell = 'D8'; Value = '0.8059' }
    @{ Cell = 'G8'; Value = '12' }
    @{ Cell = 'D9'; Value = '0.9114' }
    @{ Cell = 'G9'; Value = '12' }
    @{ Cell = 'D10'; Value = '0.1409' }
    @{ Cell = 'G10'; Value = '12' }
    @{ Cell = 'D11'; Value = '0.07444' }
    @{ Cell = 'G11'; Value = '12' }
    @{ Cell = 'D12'; Value = '0.03105' }
    @{ Cell = 'G12'; Value = '12' }
    @{ Cell = 'D13'; Value = '0.03038' }
    @{ Cell = 'G13'; Value = '12' }
    @{ Cell = 'D14'; Value = '0.09382' }
    @{ Cell = 'G14'; Value = '12' }
    @{ Cell = 'D15'; Value = '3.867' }
    @{ Cell = 'G15'; Value = '12' }
    @{ Cell = 'D16'; Value = '0.001583' }
    @{ Cell = 'G16'; Value = '12' }
    @{ Cell = 'D17'; Value = '0.04770' }
    @{ Cell = 'G17'; Value = '12' }
    @{ Cell = 'D18'; Value = '0.01828' }
    @{ Cell = 'G18'; Value = '12' }
    @{ Cell = 'D19'; Value = '0.0005814' }
    @{ Cell = 'G19'; Value = '12' }
    @{ Cell = 'D20'; Value = '0.006442' }
    @{ Cell = 'G20'; Value = '12' }
    @{ Cell = 'D21'; Value = '0.004989' }
    @{ Cell = 'G21'; Value = '12' }
    @{ Cell = 'D22'; Value = '0.0009997' }
    @{ Cell = 'G22'; Value = '12' }
    @{ Cell = 'D23'; Value = '0.0001501' }
    @{ Cell = 'G23'; Value = '12' }
    @{ Cell = 'D24'; Value = '3.699' }
    @{ Cell = 'G24'; Value = '12' }
    @{ Cell = 'D25'; Value = '2.210' }
    @{ Cell = 'G25'; Value = '12' }
    @{ Cell = 'D26'; Value = '0.3251' }
    @{ Cell = 'G26'; Value = '12' }
    @{ Cell = 'D27'; Value = '0.1307' }
    @{ Cell = 'G27'; Value = '12' }
    @{ Cell = 'G28'; Value = '12' }
    @{ Cell = 'G29'; Value = '12' }
    @{ Cell = 'G30'; Value = '12' }
    @{ Cell = 'G31'; Value = '12' }
    @{ Cell = 'G32'; Value = '12' }
    @{ Cell = 'G33'; Value = '12' }
    @{ Cell = 'G34'; Value = '12' }
    @{ Cell = 'G35'; Value = '12' }
    @{ Cell = 'G36'; Value = '12' }
    @{ Cell = 'G37'; Value = '12' }
    @{ Cell = 'G38'; Value = '12' }
    @{ Cell = 'G39'; Value = '12' }
    @{ Cell = 'D40'; Value = '0.03991' }
    @{ Cell = 'G40'; Value = '12' }
    @{ Cell = 'D41'; Value = '0.003043' }
    @{ Cell = 'E41'; Value = '40KickTokenKICKWorstin24h' }
    @{ Cell = 'G41'; Value = '12' }
    @{ Cell = 'D42'; Value = '0.1071' }
    @{ Cell = 'G42'; Value = '12' }
    @{ Cell = 'D43'; Value = '0.002682' }
    @{ Cell = 'G43'; Value = '12' }
    @{ Cell = 'D44'; Value = '0.007572' }
    @{ Cell = 'G44'; Value = '12' }
    @{ Cell = 'D45'; Value = '0.00005599' }
    @{ Cell = 'G45'; Value = '12' }
    @{ Cell = 'G46'; Value = '12' }
    @{ Cell = 'D47'; Value = '0.4993' }
    @{ Cell = 'E47'; Value = '46CoinbaseStockTokenCOIN' }
    @{ Cell = 'G47'; Value = '12' }
    @{ Cell = 'D48'; Value = '0.2098' }
    @{ Cell = 'G48'; Value = '12' }
    @{ Cell = 'D49'; Value = '0.00002101' }
    @{ Cell = 'G49'; Value = '12' }
    @{ Cell = 'D50'; Value = '0.01011' }
    @{ Cell = 'G50'; Value = '12' }
    @{ Cell = 'G51'; Value = '12' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}

Write-Output ("Updated " + $updates.Count + " cells")
